$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 17.47311853170041
$ws.Range("C2").Value = 13.31162467959333
$ws.Range("E2").Value = 17.04237244971407
$ws.Range("F2").Value = 34.80501544714441
$ws.Range("G2").Value = 25.56536607635957
$ws.Range("H2").Value = 13.58846655459367
$ws.Range("J2").Value = 7.312450490235189
$ws.Range("L2").Value = 13.02084353093219
$ws.Range("O2").Value = 20.23940658146295
$ws.Range("B3").Value = 16.79091540051647
$ws.Range("C3").Value = 13.21924130691981
$ws.Range("E3").Value = 17.08408044848404
$ws.Range("F3").Value = 34.90122112622637
$ws.Range("G3").Value = 25.73530495884855
$ws.Range("H3").Value = 13.65600376351749
$ws.Range("J3").Value = 7.306307012933496
$ws.Range("L3").Value = 12.96410139598513
$ws.Range("O3").Value = 20.3605485412556
$ws.Range("B4").Value = 16.35764425658505
$ws.Range("C4").Value = 13.16340286502873
$ws.Range("E4").Value = 17.11237135847293
$ws.Range("F4").Value = 34.97011349312481
$ws.Range("G4").Value = 25.85105282481153
$ws.Range("H4").Value = 13.70016724255202
$ws.Range("J4").Value = 7.302780660146101
$ws.Range("L4").Value = 12.93054245260712
$ws.Range("O4").Value = 20.44051829783908
$ws.Range("B5").Value = 16.1776914204773
$ws.Range("C5").Value = 13.14088763902067
$ws.Range("E5").Value = 17.12457449623178
$ws.Range("F5").Value = 35.00064937583375
$ws.Range("G5").Value = 25.90106681888133
$ws.Range("H5").Value = 13.71884200454054
$ws.Range("J5").Value = 7.301406028235746
$ws.Range("L5").Value = 12.91719759864023
$ws.Range("O5").Value = 20.4745083674243
$ws.Range("B6").Value = 16.14761277768697
$ws.Range("C6").Value = 13.13716394507194
$ws.Range("E6").Value = 17.12664154541883
$ws.Range("F6").Value = 35.00586826681618
$ws.Range("G6").Value = 25.90954283550191
$ws.Range("H6").Value = 13.72198387215628
$ws.Range("J6").Value = 7.301181556480009
$ws.Range("L6").Value = 12.91500191232572
$ws.Range("O6").Value = 20.48023695546519
$ws.Range("B7").Value = 16.35523075901868
$ws.Range("C7").Value = 13.16309822562749
$ws.Range("E7").Value = 17.11253320388836
$ws.Range("F7").Value = 34.97051535461701
$ws.Range("G7").Value = 25.85171583779434
$ws.Range("H7").Value = 13.70041635256206
$ws.Range("J7").Value = 7.302761868021694
$ws.Range("L7").Value = 12.93036112922212
$ws.Range("O7").Value = 20.44097102869805
$ws.Range("B8").Value = 17.24101153949723
$ws.Range("C8").Value = 13.27959729824658
$ws.Range("E8").Value = 17.05619668988079
$ws.Range("F8").Value = 34.83614402673771
$ws.Range("G8").Value = 25.62158080736151
$ws.Range("H8").Value = 13.61119396844592
$ws.Range("J8").Value = 7.31028156282724
$ws.Range("L8").Value = 13.0010182288
$ws.Range("O8").Value = 20.28001444406581
$ws.Range("B9").Value = 18.85509135471586
$ws.Range("C9").Value = 13.51432771559465
$ws.Range("E9").Value = 16.96700333022249
$ws.Range("F9").Value = 34.65090450485559
$ws.Range("G9").Value = 25.26174710089376
$ws.Range("H9").Value = 13.45761483168635
$ws.Range("J9").Value = 7.326952389798111
$ws.Range("L9").Value = 13.14934993433002
$ws.Range("O9").Value = 20.00887295182164
$ws.Range("B10").Value = 19.95605654745167
$ws.Range("C10").Value = 13.68959453591172
$ws.Range("E10").Value = 16.91444821501977
$ws.Range("F10").Value = 34.56292916188235
$ws.Range("G10").Value = 25.0544341957227
$ws.Range("H10").Value = 13.35781642410554
$ws.Range("J10").Value = 7.340340385426722
$ws.Range("L10").Value = 13.26374603437049
$ws.Range("O10").Value = 19.83700418889709
$ws.Range("B11").Value = 20.43681300704435
$ws.Range("C11").Value = 13.76971916511484
$ws.Range("E11").Value = 16.89335639203037
$ws.Range("F11").Value = 34.53342706078764
$ws.Range("G11").Value = 24.97278088905885
$ws.Range("H11").Value = 13.31524667421598
$ws.Range("J11").Value = 7.346671160167487
$ws.Range("L11").Value = 13.31684293051137
$ws.Range("O11").Value = 19.76480050738592
$ws.Range("B12").Value = 20.61586203034458
$ws.Range("C12").Value = 13.80009849396161
$ws.Range("E12").Value = 16.8857742506159
$ws.Range("F12").Value = 34.52377233109139
$ws.Range("G12").Value = 24.94370197737458
$ws.Range("H12").Value = 13.29953362571854
$ws.Range("J12").Value = 7.349102382199579
$ws.Range("L12").Value = 13.33709077104165
$ws.Range("O12").Value = 19.73832297131375
$ws.Range("B13").Value = 20.57743584743106
$ws.Range("C13").Value = 13.79355437819137
$ws.Range("E13").Value = 16.88738919457668
$ws.Range("F13").Value = 34.52578411478522
$ws.Range("G13").Value = 24.9498824045014
$ws.Range("H13").Value = 13.30289959547103
$ws.Range("J13").Value = 7.348577278659666
$ws.Range("L13").Value = 13.3327239317189
$ws.Range("O13").Value = 19.74398686464479
$ws.Range("B14").Value = 20.45160426711607
$ws.Range("C14").Value = 13.77221783994189
$ws.Range("E14").Value = 16.8927244907667
$ws.Range("F14").Value = 34.5326023330835
$ws.Range("G14").Value = 24.97035152248722
$ws.Range("H14").Value = 13.31394578923054
$ws.Range("J14").Value = 7.34687050295572
$ws.Range("L14").Value = 13.31850596570488
$ws.Range("O14").Value = 19.76260483419164
$ws.Range("B15").Value = 20.37413455754316
$ws.Range("C15").Value = 13.75915296560683
$ws.Range("E15").Value = 16.89604523853712
$ws.Range("F15").Value = 34.53697637622346
$ws.Range("G15").Value = 24.98312988979299
$ws.Range("H15").Value = 13.32076494211693
$ws.Range("J15").Value = 7.34582944700328
$ws.Range("L15").Value = 13.30981511302775
$ws.Range("O15").Value = 19.77412158656653
$ws.Range("B16").Value = 19.924223673593
$ws.Range("C16").Value = 13.68436453134158
$ws.Range("E16").Value = 16.91588327432901
$ws.Range("F16").Value = 34.56506925919498
$ws.Range("G16").Value = 25.06002694304818
$ws.Range("H16").Value = 13.36065541654692
$ws.Range("J16").Value = 7.339931435375373
$ws.Range("L16").Value = 13.26029638281178
$ws.Range("O16").Value = 19.84184358864076
$ws.Range("B17").Value = 19.64298760261229
$ws.Range("C17").Value = 13.63857192078421
$ws.Range("E17").Value = 16.92877444720799
$ws.Range("F17").Value = 34.58500085610193
$ws.Range("G17").Value = 25.11045755625266
$ws.Range("H17").Value = 13.38585177503349
$ws.Range("J17").Value = 7.336374309472149
$ws.Range("L17").Value = 13.23018143574018
$ws.Range("O17").Value = 19.88492369725932
$ws.Range("B18").Value = 19.47934398330528
$ws.Range("C18").Value = 13.61227152096919
$ws.Range("E18").Value = 16.93645414683212
$ws.Range("F18").Value = 34.59745480232374
$ws.Range("G18").Value = 25.14065373895271
$ws.Range("H18").Value = 13.40061031230522
$ws.Range("J18").Value = 7.33435102456867
$ws.Range("L18").Value = 13.2129603197969
$ws.Range("O18").Value = 19.91026485644993
$ws.Range("B19").Value = 19.4236172819962
$ws.Range("C19").Value = 13.60337381210397
$ws.Range("E19").Value = 16.93909988171019
$ws.Range("F19").Value = 34.6018413432566
$ws.Range("G19").Value = 25.15108132864714
$ws.Range("H19").Value = 13.40565301863202
$ws.Range("J19").Value = 7.333669891163513
$ws.Range("L19").Value = 13.20714708408815
$ws.Range("O19").Value = 19.91894143296729
$ws.Range("B20").Value = 19.67312154875993
$ws.Range("C20").Value = 13.64344279570521
$ws.Range("E20").Value = 16.92737472945782
$ws.Range("F20").Value = 34.58277662515408
$ws.Range("G20").Value = 25.10496582907671
$ws.Range("H20").Value = 13.38314202028718
$ws.Range("J20").Value = 7.336750629272879
$ws.Range("L20").Value = 13.23337693053875
$ws.Range("O20").Value = 19.88027948488014
$ws.Range("B21").Value = 20.48864639396297
$ws.Range("C21").Value = 13.77848402109677
$ws.Range("E21").Value = 16.89114639727375
$ws.Range("F21").Value = 34.53055845234891
$ws.Range("G21").Value = 24.96428910072314
$ws.Range("H21").Value = 13.3106902000915
$ws.Range("J21").Value = 7.347370910010263
$ws.Range("L21").Value = 13.32267838253709
$ws.Range("O21").Value = 19.75711278880781
$ws.Range("B22").Value = 21.00409278599244
$ws.Range("C22").Value = 13.86695328439896
$ws.Range("E22").Value = 16.86982895299083
$ws.Range("F22").Value = 34.50527499540379
$ws.Range("G22").Value = 24.88309154521794
$ws.Range("H22").Value = 13.26571237477678
$ws.Range("J22").Value = 7.354509114817698
$ws.Range("L22").Value = 13.38185987824959
$ws.Range("O22").Value = 19.68165692476382
$ws.Range("B23").Value = 20.73062887594037
$ws.Range("C23").Value = 13.81972240581521
$ws.Range("E23").Value = 16.88099056986901
$ws.Range("F23").Value = 34.51795874683955
$ws.Range("G23").Value = 24.92543823527505
$ws.Range("H23").Value = 13.28950056509256
$ws.Range("J23").Value = 7.350681504447587
$ws.Range("L23").Value = 13.35020241047426
$ws.Range("O23").Value = 19.72146641934022
$ws.Range("B24").Value = 19.65950407917194
$ws.Range("C24").Value = 13.64124058953185
$ws.Range("E24").Value = 16.92800670583407
$ws.Range("F24").Value = 34.58377910094662
$ws.Range("G24").Value = 25.10744489345534
$ws.Range("H24").Value = 13.38436625061798
$ws.Range("J24").Value = 7.336580427270937
$ws.Range("L24").Value = 13.23193195880616
$ws.Range("O24").Value = 19.88237734616289
$ws.Range("B25").Value = 18.43271693075362
$ws.Range("C25").Value = 13.45025916186201
$ws.Range("E25").Value = 16.98885379698742
$ws.Range("F25").Value = 34.69259281610578
$ws.Range("G25").Value = 25.34916434010874
$ws.Range("H25").Value = 13.49687281801639
$ws.Range("J25").Value = 7.322239834744662
$ws.Range("L25").Value = 13.10823065631968
$ws.Range("O25").Value = 20.07743869543503